$p = $ppt.ActivePresentation

# The "Wyniki testow" slide (currently at position 16) is moved up to
# position 15, ahead of the "Testy" slide - i.e. the two slides swap order.
$s = $p.Slides.Item(16)
$s.MoveTo(15)
